$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("Rules"!B11) needs to hold the TEXT string "1" (not the number 1),
# while keeping its existing cell style/format untouched. A plain
# `$ws.Range("B11").Value = "1"` would be auto-coerced into the numeric value 1
# (since "1" parses as a number), which is not what we want here.
#
# To force a genuine text value without disturbing B11's formatting/style, stage
# the text on an unused scratch cell via a formula that evaluates to the string
# "1" (TEXT() always returns text), copy that cell, and paste *values only* into
# B11. PasteSpecial(xlPasteValues) carries over the literal text value without
# touching the destination cell's existing number format/style. Finally, clear
# the scratch cell so it leaves no trace in the saved workbook.
$helper = $ws.Range("Z1")
$helper.Formula = "=TEXT(1,""0"")"

$helper.Copy()
$target = $ws.Range("B11")
$target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$helper.Clear()
$excel.CutCopyMode = $false
